$wb = $excel.ActiveWorkbook

# PowerShell COM doesn't pre-define the Excel VBA enum constants, so the
# values are spelled out numerically (same values VBA uses under the hood):
#   xlValidateList   = 3   (Type:=)
#   xlValidAlertStop = 1   (AlertStyle:=)
#   xlBetween        = 1   (Operator:=)
$xlValidateList   = 3
$xlValidAlertStop = 1
$xlBetween        = 1

# --- Sheet "sessions": add 4 new list data validations (G, H, J, L) ---
$sessions = $wb.Worksheets.Item("sessions")

$sessions.Range("G2:G1001").Validation.Add($xlValidateList, $xlValidAlertStop, $xlBetween, '"PRIVATE,RESTRICTED,SHARED,PUBLIC"')
$sessions.Range("H2:H1001").Validation.Add($xlValidateList, $xlValidAlertStop, $xlBetween, '"Lab,Home,Classroom,Outdoor,Clinic"')
$sessions.Range("J2:J1001").Validation.Add($xlValidateList, $xlValidAlertStop, $xlBetween, '"AL,AK,AZ,AR,CA,CO,CT,DE,DC,FL,GA,HI,ID,IL,IN,IA,KS,KY,LA,ME,MT,NE,NV,NH,NJ,NM,NY,NC,ND,OH,OK,OR,MD,MA,MI,MN,MS,MO,PA,RI,SC,SD,TN,TX,UT,VT,VA,WA,WV,WI,WY"')
$sessions.Range("L2:L1001").Validation.Add($xlValidateList, $xlValidAlertStop, $xlBetween, '"None,PRIVATE,SHARED,EXCERPTS,PUBLIC"')

# --- Sheet "participants": remove the "category" column (K) and shift the ---
# --- trailing "consent" column left, then add 3 new list validations ---
$participants = $wb.Worksheets.Item("participants")

# xlShiftToLeft = -4159
$participants.Range("K1").Delete(-4159)

$participants.Range("G2:G1001").Validation.Add($xlValidateList, $xlValidAlertStop, $xlBetween, '"American Indian or Alaska Native,Asian,Native Hawaiian or Other Pacific Islander,Black or African American,White,Multiple"')
$participants.Range("F2:F1001").Validation.Add($xlValidateList, $xlValidAlertStop, $xlBetween, '"Female,Male"')
$participants.Range("H2:H1001").Validation.Add($xlValidateList, $xlValidAlertStop, $xlBetween, '"Not Hispanic or Latino,Hispanic or Latino"')
